$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 48,10
$arr[0,0] = "2014-10"
$arr[0,1] = 97.79470000000001
$arr[0,2] = 100.201
$arr[0,3] = 97.17870000000001
$arr[0,4] = 98.6681
$arr[0,5] = 96.89190000000001
$arr[0,6] = 101.2659
$arr[0,7] = 100.6979
$arr[0,8] = 99.4149
$arr[0,9] = 97.0462
$arr[1,0] = "2014-11"
$arr[1,1] = 97.9622
$arr[1,2] = 100.0486
$arr[1,3] = 97.5288
$arr[1,4] = 98.52630000000001
$arr[1,5] = 96.2859
$arr[1,6] = 101.5308
$arr[1,7] = 100.9224
$arr[1,8] = 99.4944
$arr[1,9] = 97.1726
$arr[2,0] = "2014-12"
$arr[2,1] = 97.9765
$arr[2,2] = 100.0954
$arr[2,3] = 97.3125
$arr[2,4] = 98.095
$arr[2,5] = 95.9815
$arr[2,6] = 101.3617
$arr[2,7] = 101.1735
$arr[2,8] = 99.2079
$arr[2,9] = 97.1784
$arr[3,0] = "2014-01"
$arr[3,1] = 97.3078
$arr[3,2] = 98.81529999999999
$arr[3,3] = 96.56999999999999
$arr[3,4] = 98.496
$arr[3,5] = 97.6734
$arr[3,6] = 100.1907
$arr[3,7] = 100.3282
$arr[3,8] = 98.8121
$arr[3,9] = 98.0051
$arr[4,0] = "2014-02"
$arr[4,1] = 97.3164
$arr[4,2] = 99.3806
$arr[4,3] = 96.95189999999999
$arr[4,4] = 98.1544
$arr[4,5] = 97.4999
$arr[4,6] = 99.83069999999999
$arr[4,7] = 100.6109
$arr[4,8] = 98.2978
$arr[4,9] = 97.6892
$arr[5,0] = "2014-03"
$arr[5,1] = 97.57389999999999
$arr[5,2] = 99.1067
$arr[5,3] = 96.9834
$arr[5,4] = 98.0064
$arr[5,5] = 96.2332
$arr[5,6] = 100.2444
$arr[5,7] = 100.3497
$arr[5,8] = 97.86190000000001
$arr[5,9] = 97.3022
$arr[6,0] = "2014-04"
$arr[6,1] = 98.08799999999999
$arr[6,2] = 99.2811
$arr[6,3] = 96.8218
$arr[6,4] = 98.1814
$arr[6,5] = 96.6301
$arr[6,6] = 100.0172
$arr[6,7] = 100.3665
$arr[6,8] = 98.2698
$arr[6,9] = 97.0001
$arr[7,0] = "2014-05"
$arr[7,1] = 98.2705
$arr[7,2] = 99.30240000000001
$arr[7,3] = 96.8366
$arr[7,4] = 98.24590000000001
$arr[7,5] = 97.4967
$arr[7,6] = 100.3843
$arr[7,7] = 100.5291
$arr[7,8] = 98.64709999999999
$arr[7,9] = 97.2384
$arr[8,0] = "2014-06"
$arr[8,1] = 98.6986
$arr[8,2] = 99.5946
$arr[8,3] = 97.6551
$arr[8,4] = 98.2445
$arr[8,5] = 97.70869999999999
$arr[8,6] = 100.8644
$arr[8,7] = 100.2483
$arr[8,8] = 98.7182
$arr[8,9] = 97.9669
$arr[9,0] = "2014-07"
$arr[9,1] = 98.7188
$arr[9,2] = 99.8351
$arr[9,3] = 97.2187
$arr[9,4] = 98.6923
$arr[9,5] = 97.7389
$arr[9,6] = 101.2006
$arr[9,7] = 100.4378
$arr[9,8] = 99.1084
$arr[9,9] = 97.8382
$arr[10,0] = "2014-08"
$arr[10,1] = 98.4179
$arr[10,2] = 99.836
$arr[10,3] = 97.00709999999999
$arr[10,4] = 98.71559999999999
$arr[10,5] = 98.0402
$arr[10,6] = 101.1749
$arr[10,7] = 100.4997
$arr[10,8] = 99.40170000000001
$arr[10,9] = 97.4871
$arr[11,0] = "2014-09"
$arr[11,1] = 97.7967
$arr[11,2] = 99.98869999999999
$arr[11,3] = 96.89060000000001
$arr[11,4] = 98.651
$arr[11,5] = 97.2988
$arr[11,6] = 101.2518
$arr[11,7] = 100.5247
$arr[11,8] = 99.5604
$arr[11,9] = 97.0681
$arr[12,0] = "2015-10"
$arr[12,1] = 97.40000000000001
$arr[12,2] = 98.3
$arr[12,3] = 96.90000000000001
$arr[12,4] = 96.3
$arr[12,5] = 89.3
$arr[12,6] = 98.8
$arr[12,7] = 99.2
$arr[12,8] = 94.59999999999999
$arr[12,9] = 95.3
$arr[13,0] = "2015-11"
$arr[13,1] = 96.7355
$arr[13,2] = 98.2038
$arr[13,3] = 97.0127
$arr[13,4] = 96.05070000000001
$arr[13,5] = 88.8524
$arr[13,6] = 98.4205
$arr[13,7] = 98.89660000000001
$arr[13,8] = 93.7959
$arr[13,9] = 94.82989999999999
$arr[14,0] = "2015-12"
$arr[14,1] = 96.5167
$arr[14,2] = 98.14709999999999
$arr[14,3] = 97.37860000000001
$arr[14,4] = 96.2234
$arr[14,5] = 88.4235
$arr[14,6] = 98.40940000000001
$arr[14,7] = 98.84220000000001
$arr[14,8] = 92.6759
$arr[14,9] = 94.53449999999999
$arr[15,0] = "2015-01"
$arr[15,1] = 97.661
$arr[15,2] = 100.0173
$arr[15,3] = 97.5125
$arr[15,4] = 98.1412
$arr[15,5] = 94.5971
$arr[15,6] = 101.2127
$arr[15,7] = 101.0233
$arr[15,8] = 99.01600000000001
$arr[15,9] = 97.48520000000001
$arr[16,0] = "2015-02"
$arr[16,1] = 97.511
$arr[16,2] = 99.36799999999999
$arr[16,3] = 96.9966
$arr[16,4] = 97.9276
$arr[16,5] = 93.6739
$arr[16,6] = 101.034
$arr[16,7] = 101.0356
$arr[16,8] = 98.63
$arr[16,9] = 97.0848
$arr[17,0] = "2015-03"
$arr[17,1] = 97.50230000000001
$arr[17,2] = 99.16930000000001
$arr[17,3] = 97.2424
$arr[17,4] = 97.81440000000001
$arr[17,5] = 93.5175
$arr[17,6] = 100.7247
$arr[17,7] = 101.1964
$arr[17,8] = 98.4341
$arr[17,9] = 96.84569999999999
$arr[18,0] = "2015-04"
$arr[18,1] = 97.50490000000001
$arr[18,2] = 99.3978
$arr[18,3] = 96.9143
$arr[18,4] = 97.5745
$arr[18,5] = 92.6405
$arr[18,6] = 100.4108
$arr[18,7] = 100.8607
$arr[18,8] = 98.372
$arr[18,9] = 96.80719999999999
$arr[19,0] = "2015-05"
$arr[19,1] = 97.5523
$arr[19,2] = 99.32429999999999
$arr[19,3] = 97.23950000000001
$arr[19,4] = 97.3809
$arr[19,5] = 91.2774
$arr[19,6] = 100.2774
$arr[19,7] = 100.5886
$arr[19,8] = 98.30070000000001
$arr[19,9] = 96.72239999999999
$arr[20,0] = "2015-06"
$arr[20,1] = 97.248
$arr[20,2] = 99.0531
$arr[20,3] = 96.9191
$arr[20,4] = 97.5399
$arr[20,5] = 91.7119
$arr[20,6] = 99.8154
$arr[20,7] = 100.399
$arr[20,8] = 97.59529999999999
$arr[20,9] = 95.88549999999999
$arr[21,0] = "2015-07"
$arr[21,1] = 97.1009
$arr[21,2] = 98.6019
$arr[21,3] = 97.15860000000001
$arr[21,4] = 97.069
$arr[21,5] = 91.3186
$arr[21,6] = 99.0348
$arr[21,7] = 100.005
$arr[21,8] = 96.96120000000001
$arr[21,9] = 95.4354
$arr[22,0] = "2015-08"
$arr[22,1] = 97.1897
$arr[22,2] = 98.4876
$arr[22,3] = 96.9746
$arr[22,4] = 96.7927
$arr[22,5] = 90.54389999999999
$arr[22,6] = 98.8259
$arr[22,7] = 99.7747
$arr[22,8] = 95.54559999999999
$arr[22,9] = 95.4542
$arr[23,0] = "2015-09"
$arr[23,1] = 97.4456
$arr[23,2] = 98.4355
$arr[23,3] = 97.0714
$arr[23,4] = 96.65300000000001
$arr[23,5] = 90.1491
$arr[23,6] = 99.1491
$arr[23,7] = 99.26860000000001
$arr[23,8] = 94.9196
$arr[23,9] = 95.29349999999999
$arr[24,0] = "2016-10"
$arr[24,1] = 100.2
$arr[24,2] = 99.8
$arr[24,3] = 100
$arr[24,4] = 100.1
$arr[24,5] = 102.1
$arr[24,6] = 100.7
$arr[24,7] = 100.9
$arr[24,8] = 102.5
$arr[24,9] = 96.59999999999999
$arr[25,0] = "2016-11"
$arr[25,1] = 102.1
$arr[25,2] = 100.3
$arr[25,3] = 100.5
$arr[25,4] = 101.1
$arr[25,5] = 106
$arr[25,6] = 102.2
$arr[25,7] = 101.4
$arr[25,8] = 105
$arr[25,9] = 97.5
$arr[26,0] = "2016-12"
$arr[26,1] = 103.1
$arr[26,2] = 101.6
$arr[26,3] = 100.5
$arr[26,4] = 102.7
$arr[26,5] = 109.6
$arr[26,6] = 102.9
$arr[26,7] = 102
$arr[26,8] = 109.7
$arr[26,9] = 100
$arr[27,0] = "2016-01"
$arr[27,1] = 96.2796
$arr[27,2] = 98.38160000000001
$arr[27,3] = 98.72669999999999
$arr[27,4] = 95.9823
$arr[27,5] = 88.3781
$arr[27,6] = 99.2741
$arr[27,7] = 99.32559999999999
$arr[27,8] = 93.3871
$arr[27,9] = 94.6755
$arr[28,0] = "2016-02"
$arr[28,1] = 96.1606
$arr[28,2] = 98.4104
$arr[28,3] = 98.7257
$arr[28,4] = 95.8494
$arr[28,5] = 89.1275
$arr[28,6] = 99.84610000000001
$arr[28,7] = 98.97880000000001
$arr[28,8] = 94.0303
$arr[28,9] = 93.8272
$arr[29,0] = "2016-03"
$arr[29,1] = 96.9567
$arr[29,2] = 98.6112
$arr[29,3] = 98.6758
$arr[29,4] = 96.68940000000001
$arr[29,5] = 90.99290000000001
$arr[29,6] = 100.2139
$arr[29,7] = 99.1255
$arr[29,8] = 95.27970000000001
$arr[29,9] = 93.3075
$arr[30,0] = "2016-04"
$arr[30,1] = 98.1925
$arr[30,2] = 98.53149999999999
$arr[30,3] = 99.1508
$arr[30,4] = 97.5121
$arr[30,5] = 94.54640000000001
$arr[30,6] = 100.2853
$arr[30,7] = 100.1223
$arr[30,8] = 96.9062
$arr[30,9] = 93.0427
$arr[31,0] = "2016-05"
$arr[31,1] = 98.90000000000001
$arr[31,2] = 98.7
$arr[31,3] = 99.5
$arr[31,4] = 98
$arr[31,5] = 97.09999999999999
$arr[31,6] = 100.1
$arr[31,7] = 100.1
$arr[31,8] = 97.3
$arr[31,9] = 93.3
$arr[32,0] = "2016-06"
$arr[32,1] = 98.3
$arr[32,2] = 98.90000000000001
$arr[32,3] = 99.59999999999999
$arr[32,4] = 98.2
$arr[32,5] = 96.40000000000001
$arr[32,6] = 100.3
$arr[32,7] = 100.2
$arr[32,8] = 98.2
$arr[32,9] = 93.59999999999999
$arr[33,0] = "2016-07"
$arr[33,1] = 98.90000000000001
$arr[33,2] = 99.40000000000001
$arr[33,3] = 99.40000000000001
$arr[33,4] = 98.40000000000001
$arr[33,5] = 97
$arr[33,6] = 100.5
$arr[33,7] = 99.90000000000001
$arr[33,8] = 99.09999999999999
$arr[33,9] = 94.7
$arr[34,0] = "2016-08"
$arr[34,1] = 99.40000000000001
$arr[34,2] = 99.5
$arr[34,3] = 99.7
$arr[34,4] = 98.8
$arr[34,5] = 98.7
$arr[34,6] = 100.2
$arr[34,7] = 100.2
$arr[34,8] = 100.8
$arr[34,9] = 95.90000000000001
$arr[35,0] = "2016-09"
$arr[35,1] = 100.1
$arr[35,2] = 99.59999999999999
$arr[35,3] = 100.3
$arr[35,4] = 99.2
$arr[35,5] = 99.59999999999999
$arr[35,6] = 100
$arr[35,7] = 100.5
$arr[35,8] = 101.8
$arr[35,9] = 96.2
$arr[36,0] = "2017-10"
$arr[36,1] = 105.8
$arr[36,2] = 104.7
$arr[36,3] = 100.4
$arr[36,4] = 107.4
$arr[36,5] = 114.5
$arr[36,6] = 102.5
$arr[36,7] = 103.2
$arr[36,8] = 108.3
$arr[36,9] = 109
$arr[37,0] = "2017-11"
$arr[37,1] = 104.9
$arr[37,2] = 104.4
$arr[37,3] = 100.4
$arr[37,4] = 107.4
$arr[37,5] = 112.2
$arr[37,6] = 102.2
$arr[37,7] = 102.9
$arr[37,8] = 107.4
$arr[37,9] = 108.6
$arr[38,0] = "2017-12"
$arr[38,1] = 105.5
$arr[38,2] = 103.4
$arr[38,3] = 100.3
$arr[38,4] = 106.5
$arr[38,5] = 111.2
$arr[38,6] = 102.2
$arr[38,7] = 103.3
$arr[38,8] = 104.5
$arr[38,9] = 106.7
$arr[39,0] = "2017-01"
$arr[39,1] = 104.3
$arr[39,2] = 102
$arr[39,3] = 100.9
$arr[39,4] = 103.7
$arr[39,5] = 111.7
$arr[39,6] = 101.8
$arr[39,7] = 102
$arr[39,8] = 109.9
$arr[39,9] = 101.9
$arr[40,0] = "2017-02"
$arr[40,1] = 104.7
$arr[40,2] = 102.7
$arr[40,3] = 101.1
$arr[40,4] = 105
$arr[40,5] = 112.3
$arr[40,6] = 101.7
$arr[40,7] = 101.8
$arr[40,8] = 110.4
$arr[40,9] = 103.8
$arr[41,0] = "2017-03"
$arr[41,1] = 104.8
$arr[41,2] = 103.5
$arr[41,3] = 101.4
$arr[41,4] = 104.8
$arr[41,5] = 112.8
$arr[41,6] = 102.3
$arr[41,7] = 102.8
$arr[41,8] = 109.6
$arr[41,9] = 106.7
$arr[42,0] = "2017-04"
$arr[42,1] = 103.9
$arr[42,2] = 103.9
$arr[42,3] = 100.8
$arr[42,4] = 104.4
$arr[42,5] = 110.5
$arr[42,6] = 102.1
$arr[42,7] = 101.8
$arr[42,8] = 108.1
$arr[42,9] = 108.3
$arr[43,0] = "2017-05"
$arr[43,1] = 102.9
$arr[43,2] = 103.9
$arr[43,3] = 100.4
$arr[43,4] = 104.5
$arr[43,5] = 109.3
$arr[43,6] = 101.7
$arr[43,7] = 101.6
$arr[43,8] = 107.8
$arr[43,9] = 109.4
$arr[44,0] = "2017-06"
$arr[44,1] = 104.1
$arr[44,2] = 104
$arr[44,3] = 100.4
$arr[44,4] = 104.7
$arr[44,5] = 108.9
$arr[44,6] = 101.4
$arr[44,7] = 101.8
$arr[44,8] = 107
$arr[44,9] = 108.9
$arr[45,0] = "2017-07"
$arr[45,1] = 104.1
$arr[45,2] = 103.9
$arr[45,3] = 100.5
$arr[45,4] = 105.4
$arr[45,5] = 110.5
$arr[45,6] = 101.6
$arr[45,7] = 101.8
$arr[45,8] = 106.9
$arr[45,9] = 108.4
$arr[46,0] = "2017-08"
$arr[46,1] = 104.9
$arr[46,2] = 104.4
$arr[46,3] = 100.8
$arr[46,4] = 105.9
$arr[46,5] = 112.5
$arr[46,6] = 102
$arr[46,7] = 102
$arr[46,8] = 107.4
$arr[46,9] = 108.3
$arr[47,0] = "2017-09"
$arr[47,1] = 105.5
$arr[47,2] = 104.5
$arr[47,3] = 100.1
$arr[47,4] = 106.6
$arr[47,5] = 113.8
$arr[47,6] = 102.1
$arr[47,7] = 103
$arr[47,8] = 107.6
$arr[47,9] = 108.4
$ws.Range("A2:J49").Value = $arr
